# journal_de_travail.xlsx — add two new work-log entries (2024-07-09) to
# the "Journal" sheet, put the "Total" cell (I2) into elapsed-time
# ([h]:mm:ss) format, move the active selection to E13, and set the
# sheet's print setup (A4 / portrait).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal")

# --- Row 12 : 09/07/2024, 10:40 -> 18:17, "Analyse et état de l'art" ---
$ws.Range("A12").Value = 45482
$ws.Range("B12").Value = 0.44444444444444442
$ws.Range("C12").Value = 0.76180555555555562
$ws.Range("E12").Value = "Analyse et état de l'art"

# --- Row 13 : 09/07/2024, 21:40 -> 23:20, "Analyse et état de l'art" ---
$ws.Range("A13").Value = 45482
$ws.Range("B13").Value = 0.90277777777777779
$ws.Range("C13").Value = 0.97222222222222221
$ws.Range("E13").Value = "Analyse et état de l'art"

# D12/D13 already carry the shared "=C-B" formula (filled down to D37) and
# I2 = SUM(D2:D37), so both recalculate automatically from the new rows.

# Show the running total in elapsed-hours format ([h]:mm:ss, built-in
# numFmtId 46) now that it can exceed 24h.
$ws.Range("I2").NumberFormat = "[h]:mm:ss"

# Print setup: A4, portrait.
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1

# Leave the selection where the author's last edit was.
$ws.Range("E13").Select()
